$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = '일단 수준 차이 빈부격차'
$ws.Range("A4").Value = '직업 상황'
$ws.Range("A5").Value = '동해 사람 생명 위협 상황'
$ws.Range("A6").Value = '자율 주행 교통사고 책임'
$ws.Range("A7").Value = '일자리 위협'
$ws.Range("A8").Value = '무분별 데이터 수집 저작권 침해'
$ws.Range("A9").Value = '데이터 수집 가짜 뉴스 확산'
$ws.Range("A10").Value = '활용 사이버 테러 증가'
$ws.Range("A11").Value = '페이크 목소리 모방 기술 피해'
$ws.Range("A12").Value = '질문 통해 일자리 감소 발생'
$ws.Range("A13").Value = '세계 지배'
$ws.Range("A15").Value = '사람 일자리 감소'
$ws.Range("A16").Value = '발전 제어'
$ws.Range("A17").Value = '범죄 악용'
$ws.Range("A18").Value = '전쟁 범죄'
$ws.Range("A19").Value = '전쟁 악용'
$ws.Range("A20").Value = '개인정보'
$ws.Range("A21").Value = '수준 차이 빈부격차 모든 직업 사람 가난 지속'
$ws.Range("A22").Value = '페이크 통해 사칭 위험성'
$ws.Range("A23").Value = '사람 악용 남용'
$ws.Range("A24").Value = '사람 사고 능력 하락'
$ws.Range("A25").Value = '인공 지능 매우 생각'
$ws.Range("A26").Value = '활용 인공 지능 우리 배제'
$ws.Range("A27").Value = '페이크 이미지 악용 불법 활용'
$ws.Range("A28").Value = '페이크 기술 더빙 기술 악용 사람 고인 명예훼손 가능성'
$ws.Range("A29").Value = '일자리 감소 이용 얼굴 합성 페이크'
$ws.Range("A30").Value = '발달 차후 사람 일자리 가능성'
$ws.Range("A31").Value = '점점 사람과 만약 자아 통제 수도 위험'
$ws.Range("A32").Value = '발전 인간 일자리'
$ws.Range("A33").Value = '매우'
$ws.Range("A34").Value = '더욱더 발전 정말 미래 실업률 증가'
$ws.Range("A35").Value = '미래 인공 지능 학습 능력 환경 적응 능력 인간 대부분 모습 세계 적업 인간 직업 찾기 더욱 위험성'
$ws.Range("A36").Value = '우리 성인 우리 직업 우리'
$ws.Range("A37").Value = '발전 우리 인간 일자리 점점 위험성'
$ws.Range("A38").Value = '대량학살 개조 순간 사람 끼리 서로 목숨 상황'
$ws.Range("A39").Value = '정말 나중 해먹 지도 정도 현재 기술 생각 계속 가다가 직업 시리 조금 생각'
$ws.Range("A40").Value = '이제 영화에서처럼 세상 장악 수도 생각'
$ws.Range("A41").Value = '발전 인간 공존 방법'
$ws.Range("A42").Value = '직접 사람 사람 점점 자신 생각'
$ws.Range("A43").Formula = "'"
$ws.Range("A43").Style = "Normal"
$ws.Range("A44").Value = '생활화 일자리 사람'
$ws.Range("A45").Value = '일자리 간다'
$ws.Range("A46").Value = '사람'
$ws.Range("A47").Value = '일자리 위협'
$ws.Range("A48").Value = '나중 우리 일자리 차지'
$ws.Range("A49").Value = '기술 점점 발달 사람 직업'
$ws.Range("A50").Value = '사람 직업 사람 취업률 또한'
$ws.Range("A51").Value = '사람 지능 돌변'
$ws.Range("A52").Value = '요즘 페이크 때문 유명인 피해 사기 위험'
$ws.Range("A53").Value = '인간 지배'
$ws.Range("A54").Value = '사람'
$ws.Range("A55").Value = '우리 대비'
$ws.Range("A56").Value = '발전 우리 위협'
$ws.Range("A57").Value = '예술 차근차근 공격'
$ws.Range("A58").Value = '페이크 피해'
$ws.Range("A59").Value = '인간 지배'
$ws.Range("A60").Value = '일자리 감소'
$ws.Range("A61").Value = '사이버 피해'
$ws.Range("A62").Value = '인간 지배'
$ws.Range("A63").Value = '우리 폐해'
$ws.Range("A64").Value = '잘못 활용 경우 상상 범죄 수도 사람 직업 사람 점점 일이'
$ws.Range("A65").Value = '뉴스 피해자'
$ws.Range("A66").Value = '대한 개정'
$ws.Range("A67").Value = '발전 사람 공격'
$ws.Range("A68").Value = '인간 조종'
$ws.Range("A69").Value = '사람 위험'
$ws.Range("A70").Value = '실업률 증가'
$ws.Range("A71").Value = '인간 대한 공격'
$ws.Range("A72").Value = '인간 사고'
$ws.Range("A73").Value = '사람'
$ws.Range("A74").Value = '이용 사이버 범죄'
$ws.Range("A75").Value = '사람 지능'
$ws.Range("A76").Value = '사람 지배'
$ws.Range("A77").Value = '사람 공격'
$ws.Range("A78").Value = '가짜 정보'
$ws.Range("A79").Value = '대신 일자리 감소'
$ws.Range("A80").Value = '대신 일자리 감소'
$ws.Range("A81").Value = '페이크 신분 악용'
$ws.Range("A82").Value = '인간 공격'
$ws.Range("A83").Value = '페이크 사기'
$ws.Range("A84").Value = '지배'
$ws.Range("A85").Value = '위험성 만점'
$ws.Range("A86").Value = '악용 우려'
$ws.Range("A87").Value = '데이터 학습'
$ws.Range("A88").Value = '윤리 문제 발생'
$ws.Range("A89").Value = '해킹 위험'
$ws.Range("A90").Value = '도덕'
$ws.Range("A91").Value = '문학 의미'
$ws.Range("A92").Value = '세상 지배'
$ws.Range("A93").Value = '자연 파괴'
$ws.Range("A94").Value = '범죄 악용'
$ws.Range("A95").Value = '사기 발생'
$ws.Range("A96").Value = '사람 범죄 악용'
$ws.Range("A97").Value = '범죄'
$ws.Range("A98").Value = '사이버 범죄'
$ws.Range("A99").Value = '반란 일어난다'
$ws.Range("A100").Value = '오류'
$ws.Range("A101").Value = '대한 규제'
$ws.Range("A102").Value = '직업'
